# Rename headers on "First Fleet Maps" sheet to Dublin-Core-ish names,
# and rename the itemid header on the "People" sheet to "ID".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("First Fleet Maps")
$ws1.Range("A1").Value = "dc:type"
$ws1.Range("B1").Value = "ID"
$ws1.Range("C1").Value = "dc:title"
$ws1.Range("D1").Value = "caption"
$ws1.Range("E1").Value = ">dc:creator"

# Move the active selection on sheet 1 to E2 (matches the recorded UI state).
$ws1.Range("E2").Select()

$ws2 = $wb.Worksheets.Item("People")
$ws2.Range("B1").Value = "ID"
